$d = $word.ActiveDocument

# --- Paragraph 6: PT "Programa resumido" summary list (was the PT objectives text) ---
$d.Paragraphs(6).Range.Text = "1 – Introdução aos Sistemas Produtivos2 – Papel Estratégico da Produção3 – Estratégia de Produção4 – Projeto em Gestão de Produção5 – Projeto de Produtos e Serviços6 – Projeto da Rede de Operações Produtivas7 – Arranjo Físico e Fluxo"

# --- Paragraph 7: EN "Programa resumido" summary list, italic (was the EN objectives text) ---
$d.Paragraphs(7).Range.Text = "1 - Introduction to Productive Systems2 - Strategic Role of Production3 - Production Strategy4 - Project in Production Management5 - Product and Service Project6 - Production Operations Network Project7 - Physical Arrangement and Flow"

# --- Paragraph 9: ListBullet under "Docente(s) Responsavel(eis)" now holds the PT objectives text ---
$d.Paragraphs(9).Range.Text = "Introduzir os alunos nos conceitos técnicos fundamentais de um curso de Engenharia de Produção, tendo em vista a sua formação generalista voltada para os mais diversos tipos de sistemas de produção."

# --- Paragraph 11: PT "Programa" detailed list (was the PT summary list) ---
$d.Paragraphs(11).Range.Text = "1 – Introdução aos Sistemas Produtivos; Produção na Organização. Inputs, Processos de Transformação e Outputs. Tipos de Operações de Produção. Atividades da administração da produção.2 – Papel Estratégico da Produção; Papel da função produção. Objetivos de desempenho.3 – Tipos de Manufatura; Tipos básicos de Manufatura.4 – Arranjo Físico e Fluxo; Procedimento de Arranjo Físico. Tipos básicos de arranjo físico. Projeto de arranjo físico.5 – Organização do Trabalho e Métodos;Técnicas de organização e métodos de trabalho6 - Introdução ao Planejamento e Controle de Produção.Conceituação do PCP; conciliação de suprimento e demanda; natureza do suprimento e da demanda; atividades de PCP; efeito volume-variedade no PCP.7 - Introdução à qualidade e a tecnologia de processo Importância; visões; princípios de administração da qualidade total."

# --- Paragraph 12: italic paragraph now holds the EN objectives text ---
$d.Paragraphs(12).Range.Text = "Introduce students to the fundamental technical concepts of a Industrial Engineering course, with a view to their general training aimed at the most diverse types of production systems."

# --- Paragraph 14: plain paragraph under "Programa" now holds the evaluation method text ---
$d.Paragraphs(14).Range.Text = "Aulas Expositivas; trabalhos e seminários."

# --- Paragraph 19: plain paragraph under "Bibliografia" now holds the docente bullet text ---
$d.Paragraphs(19).Range.Text = "5840535 - Messias Borges Silva"

# --- Paragraph 17 ("Avaliacao" bullets: Metodo/Criterio/Norma) ---
# Each bold label keeps its text; only the value run after it changes.
# Processed in reverse chain order (Norma, then Criterio, then Metodo) so that
# no not-yet-updated "old" search text is accidentally re-matched after an
# earlier step has already inserted it as a new value elsewhere in the paragraph.
$p17 = $d.Paragraphs(17).Range
$p17.Find.Execute("NF = (MF + PR)/2, onde PR é uma prova de recuperação.", $true, $false, $false, $false, $false, $true, 0, $false, "SLACK, N. et al. Administração da produção. São Paulo: Atlas, 2002. Textos complementares serão usados durante o curso.", 2) | Out-Null
$p17 = $d.Paragraphs(17).Range
$p17.Find.Execute("MF = (0,30*P1 + 0,30*P2 + 0,40*TRAB), onde P1 e P2 são provas e TRAB é a nota média de trabalhos e seminários.", $true, $false, $false, $false, $false, $true, 0, $false, "NF = (MF + PR)/2, onde PR é uma prova de recuperação.", 2) | Out-Null
$p17 = $d.Paragraphs(17).Range
$p17.Find.Execute("Aulas Expositivas; trabalhos e seminários.", $true, $false, $false, $false, $false, $true, 0, $false, "MF = (0,30*P1 + 0,30*P2 + 0,40*TRAB), onde P1 e P2 são provas e TRAB é a nota média de trabalhos e seminários.", 2) | Out-Null

Write-Host "Edit complete"
